# Updated 'Sheet1' via CrewAI tool
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells K1:P1 ---
$ws.Range("K1").Value = "Date"
$ws.Range("L1").Value = "Name"
$ws.Range("M1").Value = "Source"
$ws.Range("N1").Value = "Field1"
$ws.Range("O1").Value = "Field2"
$ws.Range("P1").Value = "Field3"

# Match the existing bold/centered/bordered header formatting used by A1:J1
$ws.Range("A1").Copy()
$ws.Range("K1:P1").PasteSpecial(-4122) # xlPasteFormats

# --- Data rows 2-8: the "nan" value moves from column D to column F; ---
# --- columns D and E become blank; new columns K-P are blank placeholders. ---
$dataRows = 2..8
foreach ($r in $dataRows) {
    $ws.Cells.Item($r, 6).Value = "nan"                          # F = "nan" (was in D)
    $ws.Cells.Item($r, 4).ClearContents()                        # D cleared
    $ws.Cells.Item($r, 5).ClearContents()                        # E cleared
    foreach ($col in 11..16) {                                   # K..P blank placeholders
        $ws.Cells.Item($r, $col).ClearContents()
    }
}

# --- Row 9: fix casing of Entity ID, clear old H/I/J "NA" values, ---
# --- and populate the newly added K-P columns with the shifted/duplicated data. ---
$ws.Range("F9").Value = "123abx007"
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("J9").ClearContents()

$ws.Range("K9").Value = "2025-10-16 00:00:00"
$ws.Range("L9").Value = "YYY"
$ws.Range("M9").Value = "FISB"
$ws.Range("N9").Value = "NA"
$ws.Range("O9").Value = "NA"
$ws.Range("P9").Value = "NA"
